# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell that showed the
# old "Ready for handoff" status is updated to "In Translation", and the
# Status column(s) are narrowed to fit the new (shorter) text, mirroring an
# autofit/regeneration pass of the report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-locale status lives in columns E (zh-cn) and F (de-de)
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

# Narrow the now-shorter status columns to match the regenerated content.
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status lives in column C
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"
$ws2.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status lives in column C
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"
$ws3.Columns.Item(3).ColumnWidth = 12.5
